$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Summary": update the borrower name and the headline figures
# ---------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")

$wsSummary.Range("B3").Value = "Waleed Al Nuaimi"   # Name
$wsSummary.Range("B4").Value = 3512.98              # Monthly Income (AED)
$wsSummary.Range("B6").Value = 684134               # Total Assets (AED)
$wsSummary.Range("B7").Value = 186988               # Total Liabilities (AED)
$wsSummary.Range("B8").Value = 497146               # Net Worth (AED)
$wsSummary.Range("B9").Value = 3.66                 # Asset/Liability Ratio

# ---------------------------------------------------------------------
# Sheet "Assets": insert two new vehicle rows above the existing
# "Liquid Assets" row and refresh the totals
# ---------------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("Assets")

# Push the current row 2 ("Liquid Assets") and row 3 ("TOTAL ASSETS")
# down to rows 4 and 5, inserting two fresh rows at 2 and 3.
$wsAssets.Range("A2:A3").EntireRow.Insert()

# New rows inherit the formatting of the row directly below them by
# default; instead copy the formatting that row 2 (now row 4) already
# carries so the new rows look like the other data rows.
$wsAssets.Range("A4:C4").Copy()
$wsAssets.Range("A2:C2").PasteSpecial(-4122)
$wsAssets.Range("A3:C3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsAssets.Cells.Item(2, 1).Value = "Vehicles"
$wsAssets.Cells.Item(2, 2).Value = "Luxury Car"
$wsAssets.Cells.Item(2, 3).Value = 422189

$wsAssets.Cells.Item(3, 1).Value = "Vehicles"
$wsAssets.Cells.Item(3, 2).Value = "Premium Car"
$wsAssets.Cells.Item(3, 3).Value = 255024

# Former row 2 ("Liquid Assets" / "Savings Account") is now row 4
$wsAssets.Cells.Item(4, 3).Value = 6921

# Former row 3 ("TOTAL ASSETS") is now row 5
$wsAssets.Cells.Item(5, 3).Value = 684134

# ---------------------------------------------------------------------
# Sheet "Liabilities": insert a new auto-loan row above the existing
# "Credit Cards" row and refresh the totals
# ---------------------------------------------------------------------
$wsLiabilities = $wb.Worksheets.Item("Liabilities")

# Push the current row 2 ("Credit Cards") and row 3 ("TOTAL LIABILITIES")
# down to rows 3 and 4, inserting a fresh row at 2.
$wsLiabilities.Range("A2").EntireRow.Insert()

# Copy formatting from row 3 (the shifted-down "Credit Cards" row) onto
# the new row 2 so it matches the other data rows.
$wsLiabilities.Range("A3:E3").Copy()
$wsLiabilities.Range("A2:E2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsLiabilities.Cells.Item(2, 1).Value = "Auto Loans"
$wsLiabilities.Cells.Item(2, 2).Value = "Vehicle Loan 2"
$wsLiabilities.Cells.Item(2, 3).Value = 153014
$wsLiabilities.Cells.Item(2, 4).Value = 4250
$wsLiabilities.Cells.Item(2, 5).Value = 3

# Former row 2 ("Credit Cards" / "Credit Card Balance") is now row 3
$wsLiabilities.Cells.Item(3, 3).Value = 33974
$wsLiabilities.Cells.Item(3, 4).Value = 1699

# Former row 3 ("TOTAL LIABILITIES") is now row 4
$wsLiabilities.Cells.Item(4, 3).Value = 186988
